# aggiornamento fino a 1/09/2021
# Appends the new daily COVID data rows (358-366, dates 2021-08-24 .. 2021-09-01)
# to the end of the existing table on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 358; Date = 44432; B = 2; C = 8; D = 53.02578378736661 },
    @{ Row = 359; Date = 44433; B = 1; C = 8; D = 53.02578378736661 },
    @{ Row = 360; Date = 44434; B = 0; C = 7; D = 46.39756081394578 },
    @{ Row = 361; Date = 44435; B = 2; C = 7; D = 46.39756081394578 },
    @{ Row = 362; Date = 44436; B = 0; C = 6; D = 39.76933784052495 },
    @{ Row = 363; Date = 44437; B = 0; C = 5; D = 33.14111486710413 },
    @{ Row = 364; Date = 44438; B = 1; C = 6; D = 39.76933784052495 },
    @{ Row = 365; Date = 44439; B = 2; C = 6; D = 39.76933784052495 },
    @{ Row = 366; Date = 44440; B = 0; C = 5; D = 33.14111486710413 }
)

# Pick up the formatting (date style, borders, alignment) already used by the
# last populated row in column A and carry it down onto the new rows.
$lastFormattedCell = $ws.Range("A357")

foreach ($r in $newRows) {
    $rowNum = $r.Row

    $cellA = $ws.Range("A$rowNum")
    $lastFormattedCell.Copy()
    $cellA.PasteSpecial(-4122)
    $cellA.Value = $r.Date

    $ws.Range("B$rowNum").Value = $r.B
    $ws.Range("C$rowNum").Value = $r.C
    $ws.Range("D$rowNum").Value = $r.D
}
